# Renamed few transcripts. Updated the DataSheet
# Column D ("Speaker") values "RBD" -> "T" and "Student" -> "S" for the
# specific rows that held those values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rbdToT = @(2,3,5,12,14,15,16,18,19,20,22,26,27,29,30,31,33,35,36,39,40,43,45,47,54,57,58,59,60,61,63)
foreach ($r in $rbdToT) {
    $ws.Range("D$r").Value = "T"
}

$studentToS = @(32,34,38)
foreach ($r in $studentToS) {
    $ws.Range("D$r").Value = "S"
}
